$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 112
$ws1.Range("F9").Value = 8744
$ws1.Range("F10").Value = 809
$ws1.Range("F13").Value = 986
$ws1.Range("F15").Value = 46
$ws1.Range("F18").Value = 263
$ws1.Range("F21").Value = 1034

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 112
$ws4.Range("F11").Value = 8744
$ws4.Range("F12").Value = 809
$ws4.Range("F15").Value = 986
$ws4.Range("F17").Value = 46
$ws4.Range("F20").Value = 263
$ws4.Range("F23").Value = 1034
